# Updated source files to version 1.0.1
#
# This script reproduces, via Excel COM-interop calls, the changes made to
# display-fw/TouchGFX/assets/texts/texts.xlsx:
#   1. Typography sheet: a new font row ("text_graph_size") is added at row 9.
#   2. Translation sheet:
#      - The "Galden: ..." row (SingleUseId20) is removed from its original
#        position (row 9) and re-appended at the end of the table, now using
#        the new "text_graph_size" font.
#      - The "45" row (SingleUseId34) is removed entirely.
#      - The "Cooling System: ..." row (SingleUseId19) switches to the new
#        "text_graph_size" font and gets a tweaked value string.
#      - The "30" row (SingleUseId45) switches to the new "text_graph_size"
#        font.
#      - The version string is bumped from "Version: 1.0.0" to
#        "Version: 1.0.1".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Typography sheet: add the new "text_graph_size" font definition row.
# ---------------------------------------------------------------------
$typo = $wb.Worksheets.Item("Typography")

$typo.Range("B9").Value = "text_graph_size"
$typo.Range("C9").Value = "arial.ttf"
$typo.Range("D9").Value = 16
$typo.Range("E9").Value = 4
$typo.Range("F9").Value = "?"
$typo.Range("G9").Value = "-., 0123456789CF"
$typo.Range("H9").Value = ""
$typo.Range("I9").Value = ""
$typo.Range("J9").Value = ""

# ---------------------------------------------------------------------
# 2. Translation sheet edits.
# ---------------------------------------------------------------------
$trans = $wb.Worksheets.Item("Translation")

# Capture the full row of the "Galden: ..." entry (SingleUseId20) before
# removing it, so it can be re-appended at the bottom of the table.
$galdenRow = $trans.Range("B1:F100").Find("SingleUseId20").Row
$galdenB = $trans.Cells.Item($galdenRow, 2).Value2
$galdenD = $trans.Cells.Item($galdenRow, 4).Value2
$galdenE = $trans.Cells.Item($galdenRow, 5).Value2
$galdenF = $trans.Cells.Item($galdenRow, 6).Value2

# Remove that row entirely (shifts everything below it up by one).
$trans.Rows.Item($galdenRow).Delete()

# Remove the "45" entry (SingleUseId34) entirely (shifts rows below up).
$obsoleteRow = $trans.Range("B1:F100").Find("SingleUseId34").Row
$trans.Rows.Item($obsoleteRow).Delete()

# Update the "Cooling System: ..." row (SingleUseId19): new font + new text.
$coolingRow = $trans.Range("B1:F100").Find("SingleUseId19").Row
$trans.Cells.Item($coolingRow, 3).Value = "text_graph_size"
$trans.Cells.Item($coolingRow, 6).Value = "Cooling System: <value> °<value>     "

# Update the "30" row (SingleUseId45): new font.
$singleUse45Row = $trans.Range("B1:F100").Find("SingleUseId45").Row
$trans.Cells.Item($singleUse45Row, 3).Value = "text_graph_size"

# Bump the version string.
$versionRow = $trans.Range("B1:F100").Find("SingleUseId47").Row
$trans.Cells.Item($versionRow, 6).Value = "Version: 1.0.1"

# Re-append the "Galden: ..." row at the end of the table, now using the new
# "text_graph_size" font; all other fields keep their original values.
$lastRow = $trans.Cells.Item($trans.Rows.Count, 2).End(-4162).Row
$newRow = $lastRow + 1
$trans.Cells.Item($newRow, 2).Value = $galdenB
$trans.Cells.Item($newRow, 3).Value = "text_graph_size"
$trans.Cells.Item($newRow, 4).Value = $galdenD
$trans.Cells.Item($newRow, 5).Value = $galdenE
$trans.Cells.Item($newRow, 6).Value = $galdenF
